# Gamification Concept.docx - Reward System table + summary paragraph rewrite
$d = $word.ActiveDocument

# --- Reward System table (Table 1) ----------------------------------------
$tbl = $d.Tables.Item(1)

# Header row
$tbl.Cell(1,1).Range.Find.Execute("Type of Question", $true, $false, $false, $false, $false, $true, 1, $false, "Difficulty of question", 2) | Out-Null
$tbl.Cell(1,3).Range.Find.Execute("Evaluation", $true, $false, $false, $false, $false, $true, 1, $false, "Type of question", 2) | Out-Null

# Row 2 (Multiple Choice / 100 / Instant ...)
$tbl.Cell(2,1).Range.Find.Execute("Multiple Choice", $true, $false, $false, $false, $false, $true, 1, $false, "Simple", 2) | Out-Null
$tbl.Cell(2,3).Range.Find.Execute("Instant ( there is only 1 type of right answer)", $true, $false, $false, $false, $false, $true, 1, $false, "Multiple Choice", 2) | Out-Null

# Row 3 (Short open ended / 300 / Slow ...)
$tbl.Cell(3,1).Range.Find.Execute("Short open ended", $true, $false, $false, $false, $false, $true, 1, $false, "Intermediate", 2) | Out-Null
$tbl.Cell(3,2).Range.Find.Execute("300", $true, $false, $false, $false, $false, $true, 1, $false, "200", 2) | Out-Null
$tbl.Cell(3,3).Range.Find.Execute("Slow (manual reviewers are needed to check for syntax or wording error)", $true, $false, $false, $false, $false, $true, 1, $false, "Multiple Choice", 2) | Out-Null

# Row 4 (One word open ended / 200 / Normal ...)
$tbl.Cell(4,1).Range.Find.Execute("One word open ended", $true, $false, $false, $false, $false, $true, 1, $false, "Advanced", 2) | Out-Null
$tbl.Cell(4,2).Range.Find.Execute("200", $true, $false, $false, $false, $false, $true, 1, $false, "300", 2) | Out-Null
$tbl.Cell(4,3).Range.Find.Execute("Normal (one word answers need syntax checks)", $true, $false, $false, $false, $false, $true, 1, $false, "Multiple Choice", 2) | Out-Null

# --- Paragraph below the table: rewritten explanation ----------------------
$oldSummary = "With the types of questions stated with the points allocated players will save up these points and use them on a gacha machine. When they get the spin again outcome players will have to use more points to spin again. Hence players will have to keep doing the daily quiz in order to earn more points to get more spins on the gacha machine."
$newSummary = "We felt that players would have a short daily quiz where they just have to select  couple of options before completing it. Depending on the difficulty of the questions players will be rewarded more points for getting the more challenging questions correct. Since we have multiple choice as our main form of questions players won" + [char]0x2019 + "t have to take a long time to complete it. Attracting more players to take part in the quiz as it does not take as long."

$d.Content.Find.Execute($oldSummary, $true, $false, $false, $false, $false, $true, 1, $false, $newSummary, 2) | Out-Null

# --- Picture paragraph: mark the run containing the drawing as NoProof -----
# The paragraph right after the rewritten summary holds the inline drawing.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.InlineShapes.Count -gt 0) {
        $para.Range.NoProofing = 1
    }
}
